$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new blank column before column N ("Late"/"Outstanding" split),
# shifting "Late" -> O and "Outstanding" -> Q.
$ws.Columns("N:N").Select()
$ws.Columns("N:N").EntireColumn.Insert()

# The newly inserted column picks up a plain (non bestFit) width matching
# its neighbour.
$ws.Columns("N:N").ColumnWidth = 9.1666667

# Leave the cursor where the author left off after the insert.
$ws.Range("S6").Select()
